$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85

# Row 6
$ws.Range("N6").Value = 10

# Row 7
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 4
$ws.Range("S7").Value = 1.3

# Row 8
$ws.Range("S8").Value = 1.37

# Row 11
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 3.4
$ws.Range("J11").Value = 2.63
$ws.Range("X11").Value = 10
$ws.Range("AH11").Value = 12
$ws.Range("AI11").Value = 19
$ws.Range("AK11").Value = 41
$ws.Range("AQ11").Value = 34
$ws.Range("AX11").Value = 19
